# Actualización automática de tasas-transfi.xlsx
#
# Applies the daily rate-update edit:
#   - Hoja1!A1 : refresh the "Binance" conversion lines in the note
#   - tasas!N10, tasas!O10 : refresh bs/dolar rate pair
#   - tasas!N12, tasas!O12 : refresh pesos/usdt rate pair

$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note with the new Binance rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.07 = 28091.87 pesos`n✅ 28091.87 pesos = 7.02 = 963.58 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- tasas: update the rate values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 141.5
$ws2.Range("O10").Value = 3975
$ws2.Range("N12").Value = 3999.9
$ws2.Range("O12").Value = 137.2
